# chore: update Sheets via scheduled runner
# Refreshes currentAveragePrice/currentAveragePriceNQ/currentAveragePriceHQ
# (and the derived Leve price/profit columns H:N) across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets with new market-board snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7403.8887
$ws.Range("J17").Value = 7745.294
$ws.Range("L17").Value = 23235.882
$ws.Range("N17").Value = -23571.882
$ws.Range("H103").Value = 1452.6923
$ws.Range("I103").Value = 380.75
$ws.Range("J103").Value = 3167.8
$ws.Range("K103").Value = 1142.25
$ws.Range("L103").Value = 9503.400000000001
$ws.Range("M103").Value = -556.25
$ws.Range("N103").Value = -10675.4
$ws.Range("H137").Value = 13516248
$ws.Range("I137").Value = 50001744
$ws.Range("J137").Value = 3101.2964
$ws.Range("K137").Value = 150005232
$ws.Range("L137").Value = 9303.889200000001
$ws.Range("M137").Value = -150002682
$ws.Range("N137").Value = -14403.8892
$ws.Range("H138").Value = 2879.7188
$ws.Range("J138").Value = 3173.7068
$ws.Range("L138").Value = 9521.1204
$ws.Range("N138").Value = -19801.1204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3133.5115
$ws.Range("I32").Value = 2033.6097
$ws.Range("J32").Value = 18165.5
$ws.Range("K32").Value = 2033.6097
$ws.Range("L32").Value = 18165.5
$ws.Range("M32").Value = -1746.6097
$ws.Range("N32").Value = -18739.5
$ws.Range("H102").Value = 2054.6667
$ws.Range("J102").Value = 6999
$ws.Range("L102").Value = 6999
$ws.Range("N102").Value = -10243
$ws.Range("H132").Value = 2394.1592
$ws.Range("J132").Value = 3501
$ws.Range("L132").Value = 10503
$ws.Range("N132").Value = -15563

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3160.7942
$ws.Range("I86").Value = 2732.8
$ws.Range("K86").Value = 2732.8
$ws.Range("M86").Value = -1609.8
$ws.Range("H89").Value = 3160.7942
$ws.Range("I89").Value = 2732.8
$ws.Range("K89").Value = 13664
$ws.Range("M89").Value = -8048
$ws.Range("H105").Value = 7861.227
$ws.Range("I105").Value = 1976.5454
$ws.Range("J105").Value = 13745.909
$ws.Range("K105").Value = 1976.5454
$ws.Range("L105").Value = 13745.909
$ws.Range("M105").Value = -229.5454
$ws.Range("N105").Value = -17239.909
$ws.Range("H137").Value = 69699.336
$ws.Range("J137").Value = 69699.336
$ws.Range("L137").Value = 69699.336
$ws.Range("N137").Value = -79899.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31696.97
$ws.Range("I31").Value = 1885.6428
$ws.Range("J31").Value = 150942.28
$ws.Range("K31").Value = 1885.6428
$ws.Range("L31").Value = 150942.28
$ws.Range("M31").Value = -1590.6428
$ws.Range("N31").Value = -151532.28
$ws.Range("H34").Value = 31696.97
$ws.Range("I34").Value = 1885.6428
$ws.Range("J34").Value = 150942.28
$ws.Range("K34").Value = 1885.6428
$ws.Range("L34").Value = 150942.28
$ws.Range("M34").Value = -1683.6428
$ws.Range("N34").Value = -151346.28
$ws.Range("H58").Value = 3422.625
$ws.Range("J58").Value = 8313.833000000001
$ws.Range("L58").Value = 8313.833000000001
$ws.Range("N58").Value = -8719.833000000001
$ws.Range("H107").Value = 1530.5385
$ws.Range("I107").Value = 1362.2858
$ws.Range("K107").Value = 1362.2858
$ws.Range("M107").Value = 557.7141999999999
$ws.Range("H132").Value = 3898.56
$ws.Range("I132").Value = 3208.7896
$ws.Range("J132").Value = 6082.8335
$ws.Range("K132").Value = 9626.3688
$ws.Range("L132").Value = 18248.5005
$ws.Range("M132").Value = -7096.3688
$ws.Range("N132").Value = -23308.5005
$ws.Range("H134").Value = 2540.0286
$ws.Range("I134").Value = 1949.1333
$ws.Range("K134").Value = 5847.3999
$ws.Range("M134").Value = -3312.3999
$ws.Range("H136").Value = 3422.625
$ws.Range("J136").Value = 8313.833000000001
$ws.Range("L136").Value = 24941.499
$ws.Range("N136").Value = -30041.499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2382247.5
$ws.Range("I5").Value = 462.9
$ws.Range("K5").Value = 1388.7
$ws.Range("M5").Value = -1276.7
$ws.Range("H21").Value = 360.75
$ws.Range("I21").Value = 177.4
$ws.Range("J21").Value = 666.3333
$ws.Range("K21").Value = 532.2
$ws.Range("L21").Value = 1998.9999
$ws.Range("M21").Value = -359.2
$ws.Range("N21").Value = -2344.9999
$ws.Range("H126").Value = 41669976
$ws.Range("I126").Value = 2385
$ws.Range("K126").Value = 7155
$ws.Range("M126").Value = -2215
$ws.Range("H132").Value = 4062.182
$ws.Range("I132").Value = 3248.9167
$ws.Range("K132").Value = 29240.2503
$ws.Range("M132").Value = -26710.2503
$ws.Range("H135").Value = 2382247.5
$ws.Range("I135").Value = 462.9
$ws.Range("K135").Value = 4166.099999999999
$ws.Range("M135").Value = -1631.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 55555
$ws.Range("J20").Value = 55555
$ws.Range("L20").Value = 55555
$ws.Range("N20").Value = -56045
$ws.Range("H70").Value = 12999.363
$ws.Range("I70").Value = 11427.857
$ws.Range("K70").Value = 11427.857
$ws.Range("M70").Value = -11157.857
$ws.Range("H73").Value = 12999.363
$ws.Range("I73").Value = 11427.857
$ws.Range("K73").Value = 11427.857
$ws.Range("M73").Value = -10491.857
$ws.Range("H105").Value = 71078
$ws.Range("J105").Value = 71078
$ws.Range("L105").Value = 71078
$ws.Range("N105").Value = -78066
$ws.Range("H122").Value = 7657
$ws.Range("I122").Value = 6962.5356
$ws.Range("J122").Value = 8737.277
$ws.Range("K122").Value = 20887.6068
$ws.Range("L122").Value = 26211.831
$ws.Range("M122").Value = -18437.6068
$ws.Range("N122").Value = -31111.831

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10828.9
$ws.Range("I7").Value = 6555.5713
$ws.Range("K7").Value = 6555.5713
$ws.Range("M7").Value = -6443.5713
$ws.Range("H16").Value = 3937.818
$ws.Range("I16").Value = 1135.4
$ws.Range("J16").Value = 6273.1665
$ws.Range("K16").Value = 1135.4
$ws.Range("L16").Value = 6273.1665
$ws.Range("M16").Value = -965.4000000000001
$ws.Range("N16").Value = -6613.1665
$ws.Range("H40").Value = 7370.4688
$ws.Range("I40").Value = 5997.952
$ws.Range("K40").Value = 5997.952
$ws.Range("M40").Value = -5861.952
$ws.Range("H126").Value = 10828.9
$ws.Range("I126").Value = 6555.5713
$ws.Range("K126").Value = 19666.7139
$ws.Range("M126").Value = -17196.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 25028344
$ws.Range("I64").Value = 40022484
$ws.Range("K64").Value = 40022484
$ws.Range("M64").Value = -40022236
$ws.Range("H67").Value = 25028344
$ws.Range("I67").Value = 40022484
$ws.Range("K67").Value = 40022484
$ws.Range("M67").Value = -40021626
$ws.Range("H74").Value = 11625
$ws.Range("J74").Value = 11625
$ws.Range("L74").Value = 11625
$ws.Range("N74").Value = -13497
$ws.Range("H77").Value = 11625
$ws.Range("J77").Value = 11625
$ws.Range("L77").Value = 34875
$ws.Range("N77").Value = -44235
$ws.Range("H96").Value = 1178.8889
$ws.Range("I96").Value = 1125.5
$ws.Range("K96").Value = 1125.5
$ws.Range("M96").Value = 247.5
$ws.Range("H126").Value = 4080.2144
$ws.Range("I126").Value = 4239.846
$ws.Range("K126").Value = 12719.538
$ws.Range("M126").Value = -10249.538
